$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the RGB LED module string from "LEDS" to "LED_RGB" (row 6, column F)
$ws.Range("F6").Value = "LED_RGB"

# Row 16 (Arduino pin D13) used to be an "Unused pin" placeholder row; it is now
# wired up to the built-in nano LED.
$ws.Range("B16").Value = "LED"
$ws.Range("C16").Value = "DIGITAL"
$ws.Range("D16").Value = "OUTPUT"
$ws.Range("E16").Value = "NO"
$ws.Range("F16").Value = "LED_BUILTIN"
$ws.Range("G16").Value = "Built-in nano LED pin"

# G17/G18 ("328P Receiver pin/..." explanations) drop the now-unused themed
# 微软雅黑 font variant in favor of the explicit-black-color one already used
# elsewhere in the sheet.
$ws.Range("G17").Font.Color = 0
$ws.Range("G18").Font.Color = 0
